{"js": "// The \"Caption\" paragraph style (linked character style: \"CaptionChar\") was\n// redefined so captions are no longer force-centered: the explicit\n// center alignment on the style's paragraph formatting is cleared so the\n// style falls back to left (its normal/default) alignment.\nconst styles = context.document.getStyles();\nconst captionStyle = styles.getByName(\"Caption\");\n\n// Drop the forced centering (<w:jc w:val=\"center\"/> in w:pPr) from the\n// \"Caption\" style definition.\ncaptionStyle.paragraphFormat.alignment = Word.Alignment.left;\n\nawait context.sync();\n", "ps1": "# The \"Caption\" paragraph style (linked character style: \"CaptionChar\") was\n# redefined so captions are no longer force-centered: the explicit center\n# alignment on the style's paragraph formatting is cleared so the style\n# falls back to left (its normal/default) alignment.\n$d = $word.ActiveDocument\n$captionStyle = $d.Styles(\"Caption\")\n\n# Drop the forced centering (<w:jc w:val=\"center\"/> in w:pPr) from the\n# \"Caption\" style definition. 0 = wdAlignParagraphLeft.\n$captionStyle.ParagraphFormat.Alignment = 0\n"}
